$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 233.32654
$ws.Range("J17").Value = 233.32654
$ws.Range("L17").Value = 699.97962
$ws.Range("N17").Value = -1035.97962

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1655.2222
$ws.Range("I40").Value = 1568
$ws.Range("J40").Value = 1862.375
$ws.Range("K40").Value = 1568
$ws.Range("L40").Value = 1862.375
$ws.Range("M40").Value = -1393
$ws.Range("N40").Value = -2212.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7314.85
$ws.Range("I51").Value = 10899.091
$ws.Range("K51").Value = 10899.091
$ws.Range("M51").Value = -10415.091

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2914.8572
$ws.Range("I106").Value = 2984
$ws.Range("J106").Value = 2500
$ws.Range("K106").Value = 2984
$ws.Range("L106").Value = 2500
$ws.Range("M106").Value = -2353
$ws.Range("N106").Value = -3762

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 40308.42
$ws.Range("I113").Value = 73229.86
$ws.Range("J113").Value = 1900.0834
$ws.Range("K113").Value = 73229.86
$ws.Range("L113").Value = 1900.0834
$ws.Range("M113").Value = -69975.86
$ws.Range("N113").Value = -8408.0834

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1012.5294
$ws.Range("J129").Value = 1020.19354
$ws.Range("L129").Value = 3060.58062
$ws.Range("N129").Value = -13060.58062

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 2531
$ws.Range("I131").Value = 200
$ws.Range("J131").Value = 4862
$ws.Range("K131").Value = 600
$ws.Range("L131").Value = 14586
$ws.Range("M131").Value = 4440
$ws.Range("N131").Value = -24666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1540.3922
$ws.Range("I137").Value = 1433.2
$ws.Range("J137").Value = 1609.5483
$ws.Range("K137").Value = 4299.6
$ws.Range("L137").Value = 4828.644899999999
$ws.Range("M137").Value = -1749.6
$ws.Range("N137").Value = -9928.644899999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3386.889
$ws.Range("I141").Value = 2911.476
$ws.Range("J141").Value = 5050.8335
$ws.Range("K141").Value = 8734.428
$ws.Range("L141").Value = 15152.5005
$ws.Range("M141").Value = -3554.428
$ws.Range("N141").Value = -25512.5005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2458.29
$ws.Range("I32").Value = 2252.023
$ws.Range("J32").Value = 3838.6924
$ws.Range("K32").Value = 2252.023
$ws.Range("L32").Value = 3838.6924
$ws.Range("M32").Value = -1965.023
$ws.Range("N32").Value = -4412.6924

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2488.9524
$ws.Range("I45").Value = 2189.5925
$ws.Range("J45").Value = 3027.8
$ws.Range("K45").Value = 2189.5925
$ws.Range("L45").Value = 3027.8
$ws.Range("M45").Value = -1812.5925
$ws.Range("N45").Value = -3781.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1676.0392
$ws.Range("I61").Value = 818.17645
$ws.Range("J61").Value = 2104.9707
$ws.Range("K61").Value = 818.17645
$ws.Range("L61").Value = 2104.9707
$ws.Range("M61").Value = -606.17645
$ws.Range("N61").Value = -2528.9707

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 48062.75
$ws.Range("J133").Value = 48062.75
$ws.Range("L133").Value = 48062.75
$ws.Range("N133").Value = -53122.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1676.0392
$ws.Range("I136").Value = 818.17645
$ws.Range("J136").Value = 2104.9707
$ws.Range("K136").Value = 2454.52935
$ws.Range("L136").Value = 6314.9121
$ws.Range("M136").Value = 95.47064999999975
$ws.Range("N136").Value = -11414.9121

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 30031.305
$ws.Range("I20").Value = 51482.95
$ws.Range("J20").Value = 3216.75
$ws.Range("K20").Value = 51482.95
$ws.Range("L20").Value = 3216.75
$ws.Range("M20").Value = -51235.95
$ws.Range("N20").Value = -3710.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2519.1428
$ws.Range("I134").Value = 2863.9473
$ws.Range("J134").Value = 1791.2222
$ws.Range("K134").Value = 8591.841899999999
$ws.Range("L134").Value = 5373.6666
$ws.Range("M134").Value = -6056.841899999999
$ws.Range("N134").Value = -10443.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14866.324
$ws.Range("I31").Value = 32421.062
$ws.Range("K31").Value = 32421.062
$ws.Range("M31").Value = -32126.062

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 14866.324
$ws.Range("I34").Value = 32421.062
$ws.Range("K34").Value = 32421.062
$ws.Range("M34").Value = -32219.062

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1186.6
$ws.Range("I94").Value = 828.25
$ws.Range("K94").Value = 828.25
$ws.Range("M94").Value = -377.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10915.533
$ws.Range("I99").Value = 4980
$ws.Range("K99").Value = 4980
$ws.Range("M99").Value = -3482

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 10915.533
$ws.Range("I126").Value = 4980
$ws.Range("K126").Value = 14940
$ws.Range("M126").Value = -12470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3348.375
$ws.Range("I132").Value = 3341.611
$ws.Range("K132").Value = 10024.833
$ws.Range("M132").Value = -7494.832999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6167.973
$ws.Range("I5").Value = 652.6
$ws.Range("J5").Value = 17658.334
$ws.Range("K5").Value = 1957.8
$ws.Range("L5").Value = 52975.00199999999
$ws.Range("M5").Value = -1845.8
$ws.Range("N5").Value = -53199.00199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2147.0615
$ws.Range("I68").Value = 1427.5186
$ws.Range("K68").Value = 4282.5558
$ws.Range("M68").Value = -3471.5558

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2147.0615
$ws.Range("I71").Value = 1427.5186
$ws.Range("K71").Value = 12847.6674
$ws.Range("M71").Value = -8791.6674

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 7272.727
$ws.Range("I75").Value = 1000
$ws.Range("J75").Value = 7571.4287
$ws.Range("K75").Value = 3000
$ws.Range("L75").Value = 22714.2861
$ws.Range("M75").Value = -2002
$ws.Range("N75").Value = -24710.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 7272.727
$ws.Range("I78").Value = 1000
$ws.Range("J78").Value = 7571.4287
$ws.Range("K78").Value = 9000
$ws.Range("L78").Value = 68142.85830000001
$ws.Range("M78").Value = -4008
$ws.Range("N78").Value = -78126.85830000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 263281.88
$ws.Range("I107").Value = 659.96
$ws.Range("J107").Value = 497765.72
$ws.Range("K107").Value = 1979.88
$ws.Range("L107").Value = 1493297.16
$ws.Range("M107").Value = -59.88000000000011
$ws.Range("N107").Value = -1497137.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1942.8
$ws.Range("I129").Value = 625
$ws.Range("J129").Value = 2272.25
$ws.Range("K129").Value = 1875
$ws.Range("L129").Value = 6816.75
$ws.Range("M129").Value = 3125
$ws.Range("N129").Value = -16816.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1390296.1
$ws.Range("J131").Value = 1551890.1
$ws.Range("L131").Value = 4655670.300000001
$ws.Range("N131").Value = -4665750.300000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 6167.973
$ws.Range("I135").Value = 652.6
$ws.Range("J135").Value = 17658.334
$ws.Range("K135").Value = 5873.400000000001
$ws.Range("L135").Value = 158925.006
$ws.Range("M135").Value = -3338.400000000001
$ws.Range("N135").Value = -163995.006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 367.2857
$ws.Range("I2").Value = 320.16666
$ws.Range("J2").Value = 650
$ws.Range("K2").Value = 320.16666
$ws.Range("L2").Value = 650
$ws.Range("M2").Value = -207.16666
$ws.Range("N2").Value = -876

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2442.611
$ws.Range("I132").Value = 1703.7273
$ws.Range("J132").Value = 3603.7144
$ws.Range("K132").Value = 5111.1819
$ws.Range("L132").Value = 10811.1432
$ws.Range("M132").Value = -2581.1819
$ws.Range("N132").Value = -15871.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 466.66666
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -105
$ws.Range("N22").Value = -1090

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 466.66666
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 400
$ws.Range("L27").Value = 500
$ws.Range("M27").Value = -293
$ws.Range("N27").Value = -714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 41069.535
$ws.Range("I107").Value = 10655.85
$ws.Range("J107").Value = 101896.9
$ws.Range("K107").Value = 31967.55
$ws.Range("L107").Value = 305690.7
$ws.Range("M107").Value = -30047.55
$ws.Range("N107").Value = -309530.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2586
$ws.Range("J126").Value = 1750
$ws.Range("L126").Value = 5250
$ws.Range("N126").Value = -10190

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1506.6842
$ws.Range("I136").Value = 807
$ws.Range("J136").Value = 2468.75
$ws.Range("K136").Value = 2421
$ws.Range("L136").Value = 7406.25
$ws.Range("M136").Value = 129
$ws.Range("N136").Value = -12506.25
